$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# New prefold data points discovered for the J1022 set
$ws.Cells.Item(22, 7).Value = 29213      # G22
$ws.Cells.Item(22, 13).Value = 0.0000049   # M22

$ws.Cells.Item(24, 7).Value = 29303      # G24
$ws.Cells.Item(24, 13).Value = 0.0000042   # M24

$ws.Cells.Item(25, 7).Value = 29296      # G25
$ws.Cells.Item(25, 13).Value = 0.0000012   # M25

# Note about a marginal / excluded point (match the header/annotation style used
# by the other notes in column F/N, e.g. N21)
$ws.Range("N21").Copy()
$ws.Range("N23").PasteSpecial(-4122)
$ws.Range("N23").Value = "I think this one was marginial and was excluded "

# Rename the pulsar ID label (corrected J1022+101 -> J1022+1001)
$ws.Range("G14").Value = "J1022+1001"

# Update the active selection to match the author's last position
$ws.Range("G15").Select()
